# ccp_memo_template.docx formatting update
#
# 1. The stray "_GoBack" bookmark that was sitting on the Heading4
#    paragraph ("Consumer Credit Panel analysis memo") is gone — Word
#    moves this auto-bookmark to wherever the last edit happened, which
#    in this revision was in the middle of the "from" placeholder text
#    (split into "fr" | "om") inside the FROM table cell.
# 2. The "First Paragraph" and "Footnote Text" paragraph styles pick up
#    explicit Georgia/11pt formatting and spacing.

$d = $word.ActiveDocument

# --- 1 & 2: relocate the "_GoBack" bookmark -------------------------------
# Find the lower-case "from" placeholder text (inside "$from$") and split
# it into "fr" / "om", dropping the bookmark between the two halves.
# Because a document can only have one bookmark with a given name, adding
# "_GoBack" here automatically removes it from its old location on the
# Heading4 paragraph.
$rng = $d.Content
$found = $rng.Find.Execute("from", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $splitPoint = $rng.Start + 2
    $bmRange = $d.Range($splitPoint, $splitPoint)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}

# --- 3: "First Paragraph" style -------------------------------------------
$firstParagraph = $d.Styles("FirstParagraph")
$firstParagraph.BaseStyle = "Normal"
$firstParagraph.NextParagraphStyle = "Normal"
$firstParagraph.ParagraphFormat.SpaceAfter = 0
$firstParagraph.ParagraphFormat.LineSpacingRule = 3   # wdLineSpaceAtLeast
$firstParagraph.ParagraphFormat.LineSpacing = 16       # 320 twips = 16 pt
$firstParagraph.Font.Name = "Georgia"
$firstParagraph.Font.Size = 11                         # w:sz 22 half-points

# --- 4: "Footnote Text" style ----------------------------------------------
$footnoteText = $d.Styles("FootnoteText")
$footnoteText.ParagraphFormat.SpaceAfter = 6           # 120 twips = 6 pt
$footnoteText.Font.Name = "Georgia"
$footnoteText.Font.Size = 11
